# Fruta / hortaliza, semanal
# Insert a new weekly price record for "Chirimoya" (Agrícola del Norte S.A. de Arica)
# at row 21, pushing the existing rows 21-26 down to 22-27.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 21:26 down to 22:27 to make room for the new record.
$ws.Rows("21:21").Insert()

# Populate the newly inserted row 21 with the new weekly record.
$ws.Range("A21").Value2 = 1
$ws.Range("B21").Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Range("C21").Value2 = "Arica y Parinacota"
$ws.Range("D21").Value2 = 45205
$ws.Range("E21").Value2 = 15
$ws.Range("F21").Value2 = "Fruta"
$ws.Range("G21").Value2 = 100107
$ws.Range("H21").Value2 = "Otros"
$ws.Range("I21").Value2 = 100107002
$ws.Range("J21").Value2 = "Chirimoya"
$ws.Range("K21").Value2 = "Cultivar IV Región"
$ws.Range("L21").Value2 = "Primera"
$ws.Range("M21").Value2 = 200
$ws.Range("N21").Value2 = 22000
$ws.Range("O21").Value2 = 23000
$ws.Range("P21").Value2 = 22500
$ws.Range("Q21").Value2 = "`$/bandeja 10 kilos"
$ws.Range("R21").Value2 = "Región de Coquimbo"
$ws.Range("S21").Value2 = 2250
$ws.Range("T21").Value2 = 10
